$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Developer field
$ws.Range("C3").Value = "Hudson Drozdowski"

# Row 7 - Test Case 1 (__init__, Attribute set to input values.)
$ws.Range("E7").Value = "None"
$ws.Range("F7").Value = 'title = "Book Name"                                                author = "Author Name"                                       genre = FICTION'
$ws.Range("G7").Value = "The instnace is inititalized correctly, no errors."

# Expected results for rows 8 & 9 entered first
$ws.Range("G8").Value = 'ValueError("Title cannot be blank.")'
$ws.Range("G9").Value = 'raise ValueError("Author cannot be blank.")'

# Row 9 inputs (Exception raised when author is blank) entered before row 8 inputs
$ws.Range("F9").Value = 'title = "Book Name"                                                                                     author = ""                                                         genre = TRUE_CRIME'

# Row 8 inputs (Exception raised when title is blank)
$ws.Range("F8").Value = 'title = ""                                                                                     author = "Author Name"                                         genre = NON_FICTION'

$ws.Range("E8").Value = "None"
$ws.Range("E9").Value = "None"

# Row 10 - Test Case 4 (Exception raised when invalid Genre)
$ws.Range("E10").Value = "None"
$ws.Range("F10").Value = 'title = "Book Name"                                                                                     author = "Author Name"                                         genre = RANDOM_GENRE'
$ws.Range("G10").Value = 'ValueError("Invalid Genre")'

# Row 11 - Test Case 5 (returns title attribute)
$ws.Range("E11").Value = 'The object is initialized correctly     title = "Book Name"                                                author = "Author Name"                                       genre = FICTION'
$ws.Range("G11").Value = '"Book Name"'

# Row 12 - Test Case 6 (returns author attribute)
$ws.Range("E12").Value = 'The object is initialized correctly     title = "Book Name"                                                author = "Author Name"                                       genre = FICTION'
$ws.Range("G12").Value = '"Author Name"'

# Row 13 - Test Case 7 (returns Genre attribute)
$ws.Range("E13").Value = 'The object is initialized correctly     title = "Book Name"                                                author = "Author Name"                                       genre = FICTION'
$ws.Range("G13").Value = "FICTION"

# Final view state: scrolled down so row 7 is at the top, selection left on I9
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("I9").Select() | Out-Null
